$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "3.62%" "3.59%"
Replace-Text " (3.26% to 3.97%)" " (3.25% to 3.93%)"
Replace-Text "5.5%" "5.15%"
Replace-Text " (5.06% to 5.94%)" " (4.75% to 5.56%)"
Replace-Text "6.5%" "6.49%"
Replace-Text " (6.02% to 6.97%)" " (6.04% to 6.95%)"
Replace-Text "10.86%" "10.43%"
Replace-Text " (10.26% to 11.46%)" " (9.87% to 11%)"
Replace-Text "8.5%" "8.55%"
Replace-Text " (7.96% to 9.03%)" " (8.03% to 9.07%)"
Replace-Text "15.91%" "15.56%"
Replace-Text " (15.2% to 16.62%)" " (14.88% to 16.23%)"
Replace-Text "9.89%" "9.96%"
Replace-Text " (9.31% to 10.47%)" " (9.4% to 10.52%)"
Replace-Text "21.12%" "20.65%"
Replace-Text " (20.31% to 21.91%)" " (19.88% to 21.4%)"
Replace-Text "10.74%" "10.86%"
Replace-Text " (10.13% to 11.35%)" " (10.28% to 11.45%)"
Replace-Text "26.08%" "25.78%"
Replace-Text " (25.2% to 26.96%)" " (24.93% to 26.62%)"
